$d = $word.ActiveDocument

# 1. "Dashboard " + "Orang Tua" (two runs) -> "Dashboard Orang Tua" (merge into first run's text)
$d.Content.Find.Execute("Dashboard Orang Tua", $true, $false, $false, $false, $false, $true, 1, $false, "Dashboard Orang Tua", 2)

# 2. "Walikelas" -> "Orang Tua" (single run text change; use track-changes so the
#    unrelated neighbouring "Dashboard " run is not coalesced into this run)
$d.TrackRevisions = $true
$d.Content.Find.Execute("Walikelas", $true, $false, $false, $false, $false, $true, 1, $false, "Orang Tua", 2)
$d.TrackRevisions = $false
$d.AcceptAllRevisions()

# 3. "Monitoring " + "Tahfiz" (two runs) -> "Monitoring Tahfiz"
$d.Content.Find.Execute("Monitoring Tahfiz", $true, $false, $false, $false, $false, $true, 1, $false, "Monitoring Tahfiz", 2)

# 4. "Monitoring " + "Mahfudhot" (two runs) -> "Monitoring Mahfudhot"
$d.Content.Find.Execute("Monitoring Mahfudhot", $true, $false, $false, $false, $false, $true, 1, $false, "Monitoring Mahfudhot", 2)
